# Append one record (row 60) to the "Sheet1" worksheet of Combined.xlsx
# (xl/worksheets/sheet2.xml), matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$row = 60

# New record values, in column order A..AJ.
$values = [ordered]@{
    "A"  = "DIFF"
    "B"  = "xyz"
    "C"  = "KLOPD34123"
    "D"  = "-"
    "E"  = "1"
    "F"  = "12.3"
    "G"  = "123"
    "H"  = "125"
    "I"  = "36x8x2.5"
    "J"  = "12.3"
    "K"  = "-"
    "L"  = "15"
    "M"  = "Gold"
    "N"  = "Alloy Steel"
    "O"  = "Screw+anchor"
    "P"  = "-"
    "Q"  = "Danpoo"
    "R"  = "999"
    "S"  = "12"
    "T"  = "12"
    "U"  = "100"
    "V"  = "14,854.08"
    "W"  = "1254"
    "X"  = "5"
    "Y"  = "124,015"
    "Z"  = "CN"
    "AA" = "15"
    "AB" = "5"
    "AC" = "3"
    "AD" = "3"
    "AE" = "Danpoo"
    "AF" = "Home & Kitchen"
    "AG" = "Large and Bulky"
    "AH" = "FBA"
    "AI" = "2026-01-30"
    "AJ" = "Danpoo"
}

# Columns whose text would otherwise be auto-detected/converted by Excel as a
# number or date (e.g. "999", "12.3", "2026-01-30"). These must be forced to
# Text format *before* the value is assigned, so they round-trip as literal
# strings (matching the source data's "numberStoredAsText" cells) instead of
# being coerced into numeric/date cells.
$forceText = @("E","F","G","H","J","L","R","S","T","U","V","W","X","Y","AA","AB","AC","AD","AI")

foreach ($col in $values.Keys) {
    $cell = $ws.Range("$col$row")
    if ($forceText -contains $col) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $values[$col]
}
